$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target account-name labels for rows 6 through 143 (row 6 in the array is row 6 on sheet)
$labels = @(
    '融出资金',
    '拆出资金',
    '交易性金融资产',
    '其中:指定以公允价值计量且其变动计入当期损益的金融资产',
    '衍生金融资产',
    '应收票据',
    '应收账款',
    '应收款项融资',
    '预付款项',
    '应收保费',
    '应收分保账款',
    '应收分保合同准备金',
    '应收利息',
    '应收股利',
    '其他应收款',
    '应收出口退税',
    '应收补贴款',
    '内部应收款',
    '买入返售金融资产',
    '以摊余成本计量的金融资产',
    '存货',
    '以公允价值计量且其变动计入其他综合收益的金融资产',
    '合同资产',
    '划分为持有待售的资产',
    '一年内到期的非流动资产',
    '代理业务资产',
    '其他流动资产',
    '流动资产其他项目',
    '流动资产平衡项目',
    '非流动资产合计',
    '发放贷款及垫款',
    '债权投资',
    '其他债权投资',
    '长期应收款',
    '长期股权投资',
    '投资性房地产',
    '固定资产',
    '在建工程',
    '工程物资',
    '其他权益工具投资',
    '其他非流动金融资产',
    '固定资产清理',
    '生产性生物资产',
    '油气资产',
    '使用权资产',
    '无形资产',
    '开发支出',
    '商誉',
    '长期待摊费用',
    '递延所得税资产',
    '其他非流动资产',
    '以摊余成本计量的金融资产(非流动)',
    '以公允价值计量且其变动计入其他综合收益的金融资产(非流动)',
    '可供出售金融资产',
    '持有至到期投资',
    '非流动资产平衡项目',
    '非流动资产其他项目',
    '资产总计',
    '资产其他项目',
    '资产平衡项目',
    '流动负债合计',
    '短期借款',
    '向中央银行借款',
    '吸收存款及同业存放',
    '拆入资金',
    '交易性金融负债',
    '其中:指定以公允价值计量且其变动计入当期损益的金融负债',
    '衍生金融负债',
    '应付票据',
    '应付账款',
    '预收款项',
    '合同负债',
    '卖出回购金融资产款',
    '应付手续费及佣金',
    '应付职工薪酬',
    '应交税费',
    '应付利息',
    '应付股利',
    '其他应付款',
    '应付分保账款',
    '内部应付款',
    '预计流动负债',
    '保险合同准备金',
    '代理买卖证券款',
    '代理承销证券款',
    '一年内的递延收益',
    '以摊余成本计量的金融负债',
    '应付短期债券',
    '划分为持有待售的负债',
    '一年内到期的非流动负债',
    '代理业务负债',
    '其他流动负债',
    '流动负债其他项目',
    '流动负债平衡项目',
    '非流动负债合计',
    '长期借款',
    '以摊余成本计量的金融负债(非流动)',
    '应付债券',
    '其中:优先股(应付债券)',
    '其中:永续债(应付债券)',
    '租赁负债',
    '长期应付款',
    '长期应付职工薪酬',
    '专项应付款',
    '预计负债',
    '递延收益',
    '递延所得税负债',
    '其他非流动负债',
    '非流动负债其他项目',
    '非流动负债平衡项目',
    '负债合计',
    '负债其他项目',
    '负债平衡项目',
    '所有者权益合计',
    '实收资本(或股本)',
    '其他权益工具',
    '其中:优先股(其他权益工具)',
    '其中:永续债(其他权益工具)',
    '其中:其他(其他权益工具)',
    '资本公积',
    '减:库存股',
    '其他综合收益',
    '盈余公积',
    '一般风险准备',
    '专项储备',
    '未确定的投资损失',
    '未分配利润',
    '拟分配现金股利',
    '外币报表折算差额',
    '归属于母公司所有者权益合计',
    '归属于母公司股东权益其他项目',
    '归属于母公司股东权益平衡项目',
    '少数股东权益',
    '股东权益其他项目',
    '股权权益平衡项目',
    '负债和所有者权益总计',
    '负债和股东权益其他项目',
    '负债及股东权益平衡项目'
)

$startRow = 6
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
